$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in the sheet (data starts at row 2, header at row 1)
$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1

# Add the new "treatment" header in column M (13th column)
$ws.Cells.Item(1, 13).Value2() = "treatment"

# Populate the new "treatment" column (M) with the same values as the
# existing "Treatment" column (E) for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $srcValue = $ws.Cells.Item($r, 5).Value2()
    $ws.Cells.Item($r, 13).Value2() = $srcValue
}
